$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used to round-trip numeric-looking strings (e.g. "1.00", "0.999")
# into the target cells as literal TEXT instead of letting Excel coerce them to
# numbers (which would silently drop formatting like trailing zeros).
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

function Set-TextValue($cellRef, $text) {
    $helper.Value = $text
    $helper.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

# Row 2
$ws.Range("D2").Value = '67.577.99'
$ws.Range("E2").Value = '  -3.39%  '

# Row 3
$ws.Range("D3").Value = '3.781.80'
$ws.Range("E3").Value = '  -3.68%  '

# Row 4
Set-TextValue "D4" '0.999'
$ws.Range("E4").Value = '  -0.15%  '

# Row 5
Set-TextValue "D5" '597.73'
$ws.Range("E5").Value = '  -1.89%  '

# Row 6
Set-TextValue "D6" '166.91'
$ws.Range("E6").Value = '  -2.03%  '

# Row 7
$ws.Range("D7").Value = '3.782.20'
$ws.Range("E7").Value = '  -3.68%  '

# Row 8
$ws.Range("E8").Value = '  -0.08%  '

# Row 9
Set-TextValue "D9" '0.524'
$ws.Range("E9").Value = '  -2.47%  '

# Row 10
$ws.Range("E10").Value = '  -4.47%  '

# Row 11
$ws.Range("E11").Value = '  +0.22%  '

# Row 12
$ws.Range("E12").Value = '  -3.46%  '

# Row 13
Set-TextValue "D13" '0.0000258'
$ws.Range("E13").Value = '  +0.36%  '

# Row 14
Set-TextValue "D14" '36.63'
$ws.Range("E14").Value = '  -4.56%  '

# Row 15
$ws.Range("D15").Value = '4.419.04'

# Row 16
$ws.Range("D16").Value = '3.789.51'
$ws.Range("E16").Value = '  -3.64%  '

# Row 17
$ws.Range("D17").Value = '67.586.99'
$ws.Range("E17").Value = '  -3.40%  '

# Row 18
Set-TextValue "D18" '18.24'
$ws.Range("E18").Value = '  -2.59%  '

# Row 19
Set-TextValue "D19" '7.32'
$ws.Range("E19").Value = '  -4.16%  '

# Row 20
$ws.Range("E20").Value = '  -1.21%  '

# Row 21
Set-TextValue "D21" '10.90'
$ws.Range("E21").Value = '  -2.26%  '

# Row 22
Set-TextValue "D22" '463.56'
$ws.Range("E22").Value = '  -6.08%  '

# Row 23
$ws.Range("E23").Value = '  -3.14%  '

# Row 24
$ws.Range("E24").Value = '  -4.76%  '

# Row 25
Set-TextValue "D25" '82.30'
$ws.Range("E25").Value = '  -4.28%  '

# Row 26
$ws.Range("E26").Value = '  -3.70%  '

# Row 27
Set-TextValue "D27" '12.01'
$ws.Range("E27").Value = '  -2.55%  '

# Row 28
Set-TextValue "D28" '0.998'
$ws.Range("E28").Value = '  -0.19%  '

# Row 29
$ws.Range("E29").Value = '  -1.67%  '

# Row 30
$ws.Range("E30").Value = '  -2.06%  '

# Row 31
$ws.Range("D31").Value = '3.931.94'
$ws.Range("E31").Value = '  -3.70%  '

# Row 32
$ws.Range("E32").Value = '  -3.65%  '

# Row 33
Set-TextValue "D33" '31.10'
$ws.Range("E33").Value = '  -3.50%  '

# Row 34
$ws.Range("E34").Value = '  -6.89%  '

# Row 35
Set-TextValue "D35" '9.40'
$ws.Range("E35").Value = '  -2.11%  '

# Row 36
$ws.Range("D36").Value = '3.749.28'
$ws.Range("E36").Value = '  -3.63%  '

# Row 37
$ws.Range("E37").Value = '  -4.63%  '

# Row 38
Set-TextValue "D38" '3.61'
$ws.Range("E38").Value = '  +9.74%  '

# Row 39
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D39" '0.139'
$ws.Range("E39").Value = '  -2.21%  '

# Row 40
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D40" '1.01'
$ws.Range("E40").Value = '  -3.96%  '

# Row 41
Set-TextValue "D41" '5.85'
$ws.Range("E41").Value = '  -4.80%  '

# Row 42
Set-TextValue "D42" '0.998'
$ws.Range("E42").Value = '  -0.19%  '

# Row 43
$ws.Range("E43").Value = '  -5.92%  '

# Row 44
Set-TextValue "D44" '1.97'
$ws.Range("E44").Value = '  -7.94%  '

# Row 45
$ws.Range("B45").Value = 'Cosmos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D45" '8.67'
$ws.Range("E45").Value = '  -0.05%  '

# Row 46
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue "D46" '1.00'
$ws.Range("E46").Value = '  -0.01%  '

# Row 47
$ws.Range("B47").Value = 'FLOKI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-TextValue "D47" '0.000293'
$ws.Range("E47").Value = '  +5.71%  '

# Row 48
Set-TextValue "D48" '414.29'
$ws.Range("E48").Value = '  -5.30%  '

# Row 49
Set-TextValue "D49" '46.65'
$ws.Range("E49").Value = '  -3.59%  '

# Row 50
Set-TextValue "D50" '141.93'
$ws.Range("E50").Value = '  -0.79%  '

# Row 51
Set-TextValue "D51" '26.06'
$ws.Range("E51").Value = '  +2.74%  '

# Clean up the helper cell so it leaves no trace in the saved workbook.
$helper.Clear()
